$d = $word.ActiveDocument

# The report text uses curly single quotes around the ANTLR token names
# (e.g. '‘end’', '‘function’', '‘begin’'). This reverts that wording back
# to the plain token identifiers (END_KEY, FUNCTION_KEY, BEGIN_KEY).
$openQuote  = [char]0x2018
$closeQuote = [char]0x2019

# --- Replace '‘begin’' -> 'BEGIN_KEY' ---
$rngBegin = $d.Content
$rngBegin.Find.ClearFormatting()
$rngBegin.Find.Execute($openQuote + "begin" + $closeQuote, $true, $false, $false, $false, $false, `
                        $true, 1, $false, "BEGIN_KEY", 2) | Out-Null

# --- Replace '‘function’' -> 'FUNCTION_KEY' ---
$rngFunction = $d.Content
$rngFunction.Find.ClearFormatting()
$rngFunction.Find.Execute($openQuote + "function" + $closeQuote, $true, $false, $false, $false, $false, `
                           $true, 1, $false, "FUNCTION_KEY", 2) | Out-Null

# --- Replace '‘end’' -> 'END_KEY' ---
$rngEnd = $d.Content
$rngEnd.Find.ClearFormatting()
$rngEnd.Find.Execute($openQuote + "end" + $closeQuote, $true, $false, $false, $false, $false, `
                      $true, 1, $false, "END_KEY", 2) | Out-Null

# The '_GoBack' bookmark used to sit at the very end of the paragraph
# (right after the trailing "..."); in the reverted text it instead
# belongs immediately after the newly-inserted "END_KEY" token (i.e.
# right before the following ')'). Locate that spot and move the
# bookmark there.
$rngTarget = $d.Content
$rngTarget.Find.ClearFormatting()
$rngTarget.Find.Execute("END_KEY)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$targetPos = $rngTarget.Start + 7   # length of "END_KEY"

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmRange = $d.Range($targetPos, $targetPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
